$d = $word.ActiveDocument

# 1. "Ngành có mã ngành, tên ngành, mô tả ngành." -> "Ngành có mã ngành, tên ngành."
$d.Content.Find.Execute("mã ngành, tên ngành, mô tả ngành.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "mã ngành, tên ngành.", 2)

# 2. "Chuyên ngành có ... mô tả chuyên ngành." -> "... mô tả chuyên ngành, mã ngành."
$d.Content.Find.Execute("tên chuyên ngành, mô tả chuyên ngành.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "tên chuyên ngành, mô tả chuyên ngành, mã ngành.", 2)

# 3. "Cơ sở đào tạo có ... địa chỉ cơ sở, giám đốc." -> "... địa chỉ cơ sở, email, hotline"
$d.Content.Find.Execute(", giám đốc.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ", email, hotline", 2)

# 5. "Lê Nhựt Anh: Ngành, chuyên ngành" -> add ", giới tính" and ", tạo bảng" as bold runs
$rng = $d.Content
$rng.Find.Execute("Ngành, chuyên ngành", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPoint = $rng.Duplicate
$insertPoint.Collapse(0)
$r1 = $d.Range($insertPoint.Start, $insertPoint.Start)
$r1.InsertAfter(", giới tính")
$r1.Font.Bold = 1
$r2pos = $r1.End
$r2 = $d.Range($r2pos, $r2pos)
$r2.InsertAfter(", tạo bảng")
$r2.Font.Bold = 1

# 4. Insert new paragraph "Giới tính mã giới tính, giới tính (ML: male, name; FM: female, nữ)"
#    after the "Chương trình đào tạo chi tiết ... kỳ học của môn." paragraph
$find2 = $d.Content
$find2.Find.Execute("kỳ học của môn.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $find2.Paragraphs(1)
$paraEnd = $para.Range.End
$newRange = $d.Range($paraEnd, $paraEnd)
$newRange.InsertParagraphAfter()
$newRange.Collapse(0)
$newRange2 = $d.Range($newRange.End, $newRange.End)
